$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.3317991383330853
$ws.Range("C2").Value = 0.3281418799904965
$ws.Range("D2").Value = 0.3311843299736723
$ws.Range("E2").Value = 0.3299191931893509
$ws.Range("F2").Value = 0.3300085576854925
$ws.Range("G2").Value = 0.3293798702853273
$ws.Range("H2").Value = 0.3303066264904571
$ws.Range("I2").Value = 0.3317440864685732
$ws.Range("J2").Value = 0.3310448112666844
$ws.Range("K2").Value = 0.3308388198477474

# Row 3
$ws.Range("B3").Value = 37.47427129418729
$ws.Range("C3").Value = 37.91424638716187
$ws.Range("D3").Value = 37.55087964552564
$ws.Range("E3").Value = 37.70956323217771
$ws.Range("F3").Value = 37.68130106224366
$ws.Range("G3").Value = 37.76457188719129
$ws.Range("H3").Value = 37.66847521599699
$ws.Range("I3").Value = 37.50143464632552
$ws.Range("J3").Value = 37.5829673505784
$ws.Range("K3").Value = 37.59777296712473

# Row 4
$ws.Range("B4").Value = 705.1463816188133
$ws.Range("C4").Value = 727.5724505579738
$ws.Range("D4").Value = 706.636090008987
$ws.Range("E4").Value = 703.2480404788504
$ws.Range("F4").Value = 705.2459165298341
$ws.Range("G4").Value = 703.8711656691285
$ws.Range("H4").Value = 740.8012084982269
$ws.Range("I4").Value = 721.3327823019223
$ws.Range("J4").Value = 718.3730409747536
$ws.Range("K4").Value = 704.129347885532

# Row 5
$ws.Range("B5").Value = 60.68046946363139
$ws.Range("C5").Value = 64.61147260934064
$ws.Range("D5").Value = 62.99908532641394
$ws.Range("E5").Value = 62.87320152331318
$ws.Range("F5").Value = 62.02537782024667
$ws.Range("G5").Value = 60.85413683864044
$ws.Range("H5").Value = 65.94348443192656
$ws.Range("I5").Value = 64.93857132628132
$ws.Range("J5").Value = 63.74238745915884
$ws.Range("K5").Value = 59.98334541576638

# Row 6
$ws.Range("B6").Value = 19210.46994245403
$ws.Range("C6").Value = 21343.36689511951
$ws.Range("D6").Value = 20273.12359736709
$ws.Range("E6").Value = 20199.99702232266
$ws.Range("F6").Value = 19840.53234993061
$ws.Range("G6").Value = 19294.54524144724
$ws.Range("H6").Value = 22132.02716524056
$ws.Range("I6").Value = 21412.17434955624
$ws.Range("J6").Value = 20795.43899483047
$ws.Range("K6").Value = 18926.70520541714

# Row 7
$ws.Range("B7").Value = 260.410853862097
$ws.Range("C7").Value = 930.0794887206432
$ws.Range("D7").Value = 335.7368355224742
$ws.Range("E7").Value = 620.1638090482015
$ws.Range("F7").Value = 572.5546387341868
$ws.Range("G7").Value = 735.8444778816803
$ws.Range("H7").Value = 891.2916760148346
$ws.Range("I7").Value = 1281.418077615964
$ws.Range("J7").Value = 459.1801648406978
$ws.Range("K7").Value = 467.5916221710787

# Row 8
$ws.Range("B8").Value = -1421.268481132906
$ws.Range("C8").Value = -533.8447684218023
$ws.Range("D8").Value = 1042.251108660735
$ws.Range("E8").Value = -796.0435568080635
$ws.Range("F8").Value = 391.1226483822767
$ws.Range("G8").Value = 85.30201537591725
$ws.Range("H8").Value = -81.26420842582071
$ws.Range("I8").Value = -629.9580042620903
$ws.Range("J8").Value = -513.1445524827415
$ws.Range("K8").Value = -88.94888253376229

# Row 9
$ws.Range("B9").Value = 1780.952302561597
$ws.Range("C9").Value = 1831.464713149611
$ws.Range("D9").Value = 1824.6944474451
$ws.Range("E9").Value = 1813.38477275682
$ws.Range("F9").Value = 1808.732278392392
$ws.Range("G9").Value = 1780.82557086961
$ws.Range("H9").Value = 1848.560715622883
$ws.Range("I9").Value = 1827.403092669844
$ws.Range("J9").Value = 1826.780587879397
$ws.Range("K9").Value = 1768.152759130111

# Row 10
$ws.Range("B10").Value = 8674.755747950019
$ws.Range("C10").Value = 11233.83137792102
$ws.Range("D10").Value = 12342.44977239636
$ws.Range("E10").Value = 6433.138709563868
$ws.Range("F10").Value = 14169.7601712149
$ws.Range("G10").Value = 9987.92781978636
$ws.Range("H10").Value = 12045.93099727815
$ws.Range("I10").Value = 8826.889055912263
$ws.Range("J10").Value = 13275.41784384815
$ws.Range("K10").Value = 7034.150562542316

# Row 11
$ws.Range("B11").Value = 25399.79402732752
$ws.Range("C11").Value = 14342.4850927106
$ws.Range("D11").Value = -1859.914225192286
$ws.Range("E11").Value = 22366.81174153778
$ws.Range("F11").Value = 9296.416685859222
$ws.Range("G11").Value = 8527.286980740982
$ws.Range("H11").Value = 4708.737733716883
$ws.Range("I11").Value = 10760.58832047255
$ws.Range("J11").Value = 8165.800525927958
$ws.Range("K11").Value = 9290.064687650261

# Row 12
$ws.Range("B12").Value = -4.05793232255743
$ws.Range("C12").Value = -4.061840138201696
$ws.Range("D12").Value = -4.03495637691328
$ws.Range("E12").Value = -4.018154340024099
$ws.Range("F12").Value = -4.001964213188653
$ws.Range("G12").Value = -4.075297948926083
$ws.Range("H12").Value = -4.066038691173847
$ws.Range("I12").Value = -4.083348102043941
$ws.Range("J12").Value = -4.036703183983004
$ws.Range("K12").Value = -4.050618792734797

# Row 13
$ws.Range("B13").Value = -1.530006850274178
$ws.Range("C13").Value = -1.792962178946633
$ws.Range("D13").Value = -1.756144144718701
$ws.Range("E13").Value = -1.549092089934141
$ws.Range("F13").Value = -1.482546792128378
$ws.Range("G13").Value = -1.171887799373307
$ws.Range("H13").Value = -1.703913565391732
$ws.Range("I13").Value = -1.547000580574477
$ws.Range("J13").Value = -1.575494952934238
$ws.Range("K13").Value = -1.815030835728432

# Row 14
$ws.Range("B14").Value = -1.475808784188639
$ws.Range("C14").Value = -1.740832774880834
$ws.Range("D14").Value = -1.70271213482388
$ws.Range("E14").Value = -1.496499124657895
$ws.Range("F14").Value = -1.429789635470771
$ws.Range("G14").Value = -1.119162800531754
$ws.Range("H14").Value = -1.650519244412558
$ws.Range("I14").Value = -1.493450466131732
$ws.Range("J14").Value = -1.522086793312159
$ws.Range("K14").Value = -1.761283318393644

# Row 15
$ws.Range("B15").Value = 1.966333247576683
$ws.Range("C15").Value = 1.576397429922822
$ws.Range("D15").Value = 1.690909787411947
$ws.Range("E15").Value = 1.846824292435456
$ws.Range("F15").Value = 1.924058244574713
$ws.Range("G15").Value = 2.233603813242816
$ws.Range("H15").Value = 1.742226122035485
$ws.Range("I15").Value = 1.906478029449788
$ws.Range("J15").Value = 1.870014231085878
$ws.Range("K15").Value = 1.652747983923168
